$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$firstRow = 23
$lastRow  = 66

# Snapshot existing rows 23..66 (whole A:T row) BEFORE any writes, since we are
# shifting every row down by one and inserting fresh data at the top.
$snapshot = @{}
for ($r = $firstRow; $r -le $lastRow; $r++) {
    $rowVals = @{}
    for ($c = 1; $c -le 20; $c++) {
        $rowVals[$c] = $ws.Cells.Item($r, $c).Value2
    }
    $snapshot[$r] = $rowVals
}

# Shift rows down: new row (r+1) gets old row r's data, working from the
# bottom up so we never overwrite a row before it has been captured (we
# already snapshotted everything above, but keep the write order safe too).
for ($r = $lastRow; $r -ge $firstRow; $r--) {
    $destRow = $r + 1
    $src = $snapshot[$r]
    for ($c = 1; $c -le 20; $c++) {
        $ws.Cells.Item($destRow, $c).Value2 = $src[$c]
    }
    # Carry the date number format from the source row's date cell (col D = 4)
    $ws.Cells.Item($destRow, 4).NumberFormat = $ws.Cells.Item($r, 4).NumberFormat
}

# Row 23 gets the brand-new record.
$ws.Range("A23").Value2 = 10
$ws.Range("B23").Value2 = "Vega Modelo de Temuco"
$ws.Range("C23").Value2 = "La Araucanía"
$ws.Range("D23").Value2 = 44533
$ws.Range("E23").Value2 = 9
$ws.Range("F23").Value2 = "Fruta"
$ws.Range("G23").Value2 = 100101
$ws.Range("H23").Value2 = "Berries"
$ws.Range("I23").Value2 = 100101001
$ws.Range("J23").Value2 = "Arándano (blue)"
$ws.Range("K23").Value2 = "Sin especificar"
$ws.Range("L23").Value2 = "Primera"
$ws.Range("M23").Value2 = 155
$ws.Range("N23").Value2 = 3000
$ws.Range("O23").Value2 = 3000
$ws.Range("P23").Value2 = 3000
$ws.Range("Q23").Value2 = "$/kilo"
$ws.Range("R23").Value2 = "Región del Maule"
$ws.Range("S23").Value2 = 3000
$ws.Range("T23").Value2 = 1
$ws.Range("D23").NumberFormat = "YYYY-MM-DD HH:MM:SS"
